# Weekly update: insert a new price record for the latest week at row 36,
# pushing all subsequent rows (old rows 36-74) down by one (to rows 37-75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36 (shifts rows 36..74 down to 37..75)
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the latest week's data
$ws.Range("A36").Value = 4
$ws.Range("B36").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C36").Value = "Los Lagos"
$ws.Range("D36").Value = 45128
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = 100112012
$ws.Range("G36").Value = "Espinaca"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 25
$ws.Range("K36").Value = 12000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 12000
$ws.Range("N36").Value = "$/cuna 10 kilos"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 1200
$ws.Range("Q36").Value = 10
$ws.Range("R36").Value = "Hortaliza"
